# Add a new quarterly date column (BB) to the yoy_rt_data export sheet.
# - BB1 gets the next date in series (45986), using the same style/number
#   format as the rest of the header row (copied from BA1).
# - BB3:BB21 repeat the last observed value from the corresponding BA cell
#   (carry-forward of the latest reading), matching column BA row-for-row.
# - Rows 2 and 22 have no data beyond column A, so no BB cell is written there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date value (serial date number) for column BB, row 1.
# Copy BA1's formatting (bold, border, centered, custom date number format)
# onto BB1 first, then set its value.
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("BB1").Value = 45986

# Carry the last BA value into the new BB column for rows 3 through 21.
for ($r = 3; $r -le 21; $r++) {
    $baCell = $ws.Cells.Item($r, 53)   # column BA = 53
    $bbCell = $ws.Cells.Item($r, 54)   # column BB = 54
    $baValue = $baCell.Value2
    $bbCell.Value = $baValue
}
